$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to stage numeric-looking text so Excel keeps it as a
# string (t="s"/"inlineStr") instead of auto-converting to a number.
$tmp = $ws.Range("Z1")
$tmp.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $tmp.Value = $text
    $tmp.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "243.95"
Set-TextValue $ws.Range("G2") "9"

Set-TextValue $ws.Range("D3") "23.79"
Set-TextValue $ws.Range("G3") "9"

Set-TextValue $ws.Range("D4") "5.149"
Set-TextValue $ws.Range("G4") "9"

Set-TextValue $ws.Range("D5") "0.05749"
Set-TextValue $ws.Range("G5") "9"

Set-TextValue $ws.Range("D6") "6.472"
Set-TextValue $ws.Range("G6") "9"

Set-TextValue $ws.Range("D7") "3.120"
Set-TextValue $ws.Range("G7") "9"

Set-TextValue $ws.Range("D8") "0.8104"
Set-TextValue $ws.Range("G8") "9"

Set-TextValue $ws.Range("D9") "0.8422"
Set-TextValue $ws.Range("G9") "9"

Set-TextValue $ws.Range("D10") "0.1336"
Set-TextValue $ws.Range("G10") "9"

Set-TextValue $ws.Range("D11") "0.06936"
Set-TextValue $ws.Range("G11") "9"

Set-TextValue $ws.Range("D12") "0.03124"
Set-TextValue $ws.Range("G12") "9"

Set-TextValue $ws.Range("D13") "0.02856"
Set-TextValue $ws.Range("G13") "9"

Set-TextValue $ws.Range("D14") "0.09363"
Set-TextValue $ws.Range("G14") "9"

Set-TextValue $ws.Range("D15") "3.756"
Set-TextValue $ws.Range("G15") "9"

Set-TextValue $ws.Range("D16") "0.001522"
Set-TextValue $ws.Range("G16") "9"

Set-TextValue $ws.Range("D17") "0.04650"
Set-TextValue $ws.Range("G17") "9"

Set-TextValue $ws.Range("D18") "0.0005969"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue $ws.Range("G18") "9"

Set-TextValue $ws.Range("D19") "0.006166"
Set-TextValue $ws.Range("G19") "9"

Set-TextValue $ws.Range("D20") "0.001239"
Set-TextValue $ws.Range("G20") "9"

Set-TextValue $ws.Range("D21") "0.004277"
Set-TextValue $ws.Range("G21") "9"

Set-TextValue $ws.Range("G22") "9"

Set-TextValue $ws.Range("D23") "3.501"
Set-TextValue $ws.Range("G23") "9"

Set-TextValue $ws.Range("D24") "2.083"
Set-TextValue $ws.Range("G24") "9"

Set-TextValue $ws.Range("D25") "0.3175"
Set-TextValue $ws.Range("G25") "9"

Set-TextValue $ws.Range("D26") "0.1337"
Set-TextValue $ws.Range("G26") "9"

Set-TextValue $ws.Range("G27") "9"

Set-TextValue $ws.Range("D28") "0.0002328"
Set-TextValue $ws.Range("G28") "9"

Set-TextValue $ws.Range("G29") "9"

Set-TextValue $ws.Range("G30") "9"

Set-TextValue $ws.Range("G31") "9"

Set-TextValue $ws.Range("G32") "9"

Set-TextValue $ws.Range("G33") "9"

Set-TextValue $ws.Range("G34") "9"

Set-TextValue $ws.Range("G35") "9"

Set-TextValue $ws.Range("G36") "9"

Set-TextValue $ws.Range("G37") "9"

Set-TextValue $ws.Range("G38") "9"

Set-TextValue $ws.Range("G39") "9"

Set-TextValue $ws.Range("D40") "0.03614"
Set-TextValue $ws.Range("G40") "9"

Set-TextValue $ws.Range("G41") "9"

Set-TextValue $ws.Range("G42") "9"

Set-TextValue $ws.Range("D43") "0.002879"
Set-TextValue $ws.Range("G43") "9"

Set-TextValue $ws.Range("D44") "0.007389"
Set-TextValue $ws.Range("G44") "9"

Set-TextValue $ws.Range("D45") "0.00005306"
Set-TextValue $ws.Range("G45") "9"

Set-TextValue $ws.Range("G46") "9"

Set-TextValue $ws.Range("D47") "0.2799"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue $ws.Range("G47") "9"

Set-TextValue $ws.Range("D48") "0.002276"
Set-TextValue $ws.Range("G48") "9"

Set-TextValue $ws.Range("G49") "9"

Set-TextValue $ws.Range("G50") "9"

Set-TextValue $ws.Range("G51") "9"

$tmp.Clear()
$excel.CutCopyMode = 0
